# Apply the "uploadPage" shipment-info column additions to 시트1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Stamp format onto the new cells first, copying the sheet's
#        current (pre-edit) header/data look from E1 before that look is
#        reassigned below. This gives the brand-new columns the same
#        "look" the original 5 columns used to have.
$newCols = @(6, 7, 8, 9, 10, 11, 12)   # F:L
foreach ($col in $newCols) {
    $ws.Range("E1").Copy()
    $ws.Cells.Item(1, $col).PasteSpecial(-4122)   # row 1 (headers)
    $ws.Range("E1").Copy()
    $ws.Cells.Item(2, $col).PasteSpecial(-4122)   # row 2 (data)
}

# --- 2) Re-color the original columns (A:E) to a plain black font, which
#        is the look the edited file now uses for the pre-existing data.
$ws.Range("A1:E1").Font.Color = 0
$ws.Range("A2:B2").Font.Color = 0
$ws.Range("D2:E2").Font.Color = 0
$ws.Range("C2").Font.Color = 0

# --- 3) Header text -------------------------------------------------------
$ws.Cells.Item(1, 5).Value = "택배회사명"          # E1 : 할말 -> 택배회사명
$ws.Cells.Item(1, 6).Value = "택배배송시간"
$ws.Cells.Item(1, 7).Value = "송장번호"
$ws.Cells.Item(1, 8).Value = "주문번호"
$ws.Cells.Item(1, 9).Value = "구/면"
$ws.Cells.Item(1, 10).Value = "동/리"
$ws.Cells.Item(1, 11).Value = "배송예정일"
$ws.Cells.Item(1, 12).Value = "결제금액"

# --- 4) Data row values ----------------------------------------------------
$ws.Cells.Item(2, 5).Value = "CJ택배"               # E2 : 보내지려나! -> CJ택배
$ws.Cells.Item(2, 7).Value = 2901248912             # G2 : 송장번호
$ws.Cells.Item(2, 8).Value = 912399                 # H2 : 주문번호
$ws.Cells.Item(2, 9).Value = "처인구"                # I2 : 구/면
$ws.Cells.Item(2, 10).Value = "왕곡동"               # J2 : 동/리
$ws.Cells.Item(2, 11).Value = "2일뒤"                # K2 : 배송예정일
$ws.Cells.Item(2, 12).Value = 33000                 # L2 : 결제금액

# F2 : delivery time, stored as a time-of-day serial value formatted am/pm h:mm.
$ws.Cells.Item(2, 6).Value = 0.5416666666666666
$ws.Cells.Item(2, 6).NumberFormat = "am/pm h:mm"
